# Word COM-interop script implementing the commit:
#   "shortened pdf to fit onto 1 page"
#
# Changes:
#  1. Merge the two runs split by a "_GoBack" bookmark in the "Im zweiten
#     Raum ..." paragraph back into a single run (removing the bookmark).
#  2. "Modul 1" -> "Angreifbar"
#  3. "Modul 2" -> "Attribute", remove the following "…." bullet
#     paragraph entirely, and add a (collapsed) "_GoBack" bookmark right
#     after the new "Attribute" run.
#  4. Remove the two blank paragraphs right before
#     "Aufgabenbereich 2: Spielen des Abenteuers".
#  5. Remove the stray <w:lastRenderedPageBreak/> before
#     "Durchführen der Kampfrunden".

$d = $word.ActiveDocument

# --- 1. Merge runs split by the _GoBack bookmark in the room-2 paragraph ---
# Replacing text that spans the bookmark collapses it back into one run and
# drops the bookmark (the concatenated text is unchanged: "öffne" + "t" =
# "öffnet").
$d.Content.Find.Execute("öffnet sich die Tür", $true, $false, $false, $false, $false,
                         $true, 1, $false, "öffnet sich die Tür", 2)

# --- 2. "Modul 1" -> "Angreifbar" ---
$d.Content.Find.Execute("Modul 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Angreifbar", 2)

# --- 3a. "Modul 2" -> "Attribute" ---
$d.Content.Find.Execute("Modul 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Attribute", 2)

# --- 3b. Remove the following "…." sub-bullet paragraph entirely ---
$findDots = $d.Content
$findDots.Find.Execute("….")
$dotsIdx = $findDots.Paragraphs.Item(1).Index
$dotsPara = $d.Paragraphs.Item($dotsIdx)
$afterPara = $d.Paragraphs.Item($dotsIdx + 1)
$d.Range($dotsPara.Range.Start, $afterPara.Range.Start).Delete()

# --- 3c. Re-locate "Attribute" and place a collapsed "_GoBack" bookmark
#         right after it (still inside its paragraph, before the
#         paragraph mark). A bookmark range collapsed exactly onto a
#         paragraph-mark position is unreliable in this engine, so we
#         temporarily insert a placeholder character after the word,
#         anchor the bookmark just before it, then remove the
#         placeholder. ---
$attrRange = $d.Content
$attrRange.Find.Execute("Attribute")
$endPos = $attrRange.End
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($endPos, $endPos + 1).Delete()

# --- 4. Remove the two blank paragraphs before "Aufgabenbereich 2" ---
$findHeading = $d.Content
$findHeading.Find.Execute("Aufgabenbereich 2: Spielen des Abenteuers")
$headingIdx = $findHeading.Paragraphs.Item(1).Index
$blank2 = $d.Paragraphs.Item($headingIdx - 1)
$blank2.Range.Delete()
$blank1 = $d.Paragraphs.Item($headingIdx - 2)
$blank1.Range.Delete()

# --- 5. Drop the stray lastRenderedPageBreak before "Durchführen der Kampfrunden" ---
$d.Content.Find.Execute("Durchführen der Kampfrunden", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Durchführen der Kampfrunden", 2)
